$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new snapshot column right before the existing "nom" / "url_produit"
# columns (currently FQ:FR). This shifts FQ->FR and FR->FS, matching the target layout.
$ws.Columns("FQ:FQ").Insert()

# Header for the freshly inserted column: the new scrape timestamp.
$ws.Range("FQ1").Value2 = "2026-02-04 22:15:34"

# The new column carries the same price snapshot as the previous column (FP),
# for every data row (rows 2 through 208). Rows whose FP cell is blank stay blank.
$lastRow = 208
for ($r = 2; $r -le $lastRow; $r++) {
    $prev = $ws.Cells.Item($r, 172)   # column FP
    $cur  = $ws.Cells.Item($r, 173)   # column FQ (new)
    $cur.Value2 = $prev.Value2
}
